# "Add files via upload" — refreshed red_ahash.xlsx results.
# The perceptual-hash column now stores the raw 64/68-bit binary string
# (instead of its signed-decimal interpretation) and the per-image timing /
# Hamming-distance columns were re-measured.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds hash strings that are all-digits but far longer than a
# double can hold exactly (64-68 characters). Pre-format the range as Text
# so Excel/COM stores the literal characters instead of rounding them into
# scientific notation, same as the original data.
$ws.Range("B2:B64").NumberFormat = "@"

$hashRows2to56 = "1110001111000001100000011000000110000001100000011100001111100011"
$hashRows57to64 = "1110001111000001100000011000000110000001100010011100001111100011"

for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 2).Value = $hashRows2to56
}
for ($r = 57; $r -le 64; $r++) {
    $ws.Cells.Item($r, 2).Value = $hashRows57to64
}

# Column C ("Время обработки" / processing time) re-measured per row.
$processingTime = @{
    2  = 0.014261
    3  = 0.001323
    4  = 0.006845
    5  = 0
    6  = 0
    7  = 0.005124
    8  = 0.001019
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0.005723
    14 = 0.00033
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0.006894
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0.007073
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0.000619
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0.006219
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0.00706
    45 = 0
    46 = 0
    47 = 0
    48 = 0.006473
    49 = 0.001216
    50 = 0
    51 = 0.006844
    52 = 0
    53 = 0.006957
    54 = 0
    55 = 0.006881
    56 = 0
    57 = 0.005925
    58 = 0.00094
    59 = 0
    60 = 0.006942
    61 = 0
    62 = 0.00708
    63 = 0.001211
    64 = 0.005541
}

foreach ($r in $processingTime.Keys) {
    $ws.Cells.Item($r, 3).Value = $processingTime[$r]
}

# Column D ("Хэммингово расстояние") — rows 57-64 go from 14 down to 1
# now that the hash is compared bit-for-bit instead of as a decimal id.
for ($r = 57; $r -le 64; $r++) {
    $ws.Cells.Item($r, 4).Value = 1
}

Write-Host "Updated hashes, processing times, and Hamming distances for rows 2-64."
